# add - management command - generate excel dummy transactions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old dummy rows (rows 2-5) before regenerating.
$ws.Range("A2:C5").ClearContents()

# Header code / product variation code / quantity dummy data.
$data = @(
    @("P-REQ-1", "PV-001", 100),
    @("P-REQ-1", "PV-002", 100),
    @("P-REQ-1", "PV-003", 100),
    @("P-REQ-2", "PV-001", 100),
    @("P-REQ-3", "PV-002", 100),
    @("P-REQ-4", "PV-003", 100),
    @("P-REQ-5", "PV-001", 100),
    @("P-REQ-6", "PV-002", 100),
    @("P-REQ-7", "PV-003", 100),
    @("P-REQ-8", "PV-001", 100),
    @("P-REQ-9", "PV-003", 100),
    @("P-REQ-10", "PV-001", 100)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

$ws.Range("A2:C13").Select()
